$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.66920033333333
$ws.Range("H2").Value = 47.00760099999999
$ws.Range("I2").Value = 0.2925937299273087
$ws.Range("J2").Value = 0.2925937299273087
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.96029466666667
$ws.Range("N2").Value = 95.88088399999999
$ws.Range("O2").Value = 0.3907265741426954
$ws.Range("P2").Value = 0.3907265741426953
$ws.Range("Q2").Value = 500.7922598443648
$ws.Range("R2").Value = 4507.130338599283
$ws.Range("S2").Value = 0.1143241457101304
$ws.Range("T2").Value = 0.1143241457101304

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.66920033333333
$ws.Range("H3").Value = 47.00760099999999
$ws.Range("I3").Value = 0.2925937299273087
$ws.Range("J3").Value = 0.2925937299273087
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 26.34807
$ws.Range("N3").Value = 79.04420999999999
$ws.Range("O3").Value = 0.3221150253382706
$ws.Range("P3").Value = 0.3221150253382706
$ws.Range("Q3").Value = 412.8531872266899
$ws.Range("R3").Value = 3715.678685040209
$ws.Range("S3").Value = 0.09424883672935414
$ws.Range("T3").Value = 0.09424883672935414

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.66920033333333
$ws.Range("H4").Value = 47.00760099999999
$ws.Range("I4").Value = 0.2925937299273087
$ws.Range("J4").Value = 0.2925937299273087
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.48872
$ws.Range("N4").Value = 70.46616
$ws.Range("O4").Value = 0.287158400519034
$ws.Range("P4").Value = 0.287158400519034
$ws.Range("Q4").Value = 368.0494592535733
$ws.Range("R4").Value = 3312.44513328216
$ws.Range("S4").Value = 0.08402074748782418
$ws.Range("T4").Value = 0.08402074748782416

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.47676966666667
$ws.Range("H5").Value = 49.43030900000001
$ws.Range("I5").Value = 0.3076736139282969
$ws.Range("J5").Value = 0.3076736139282968
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 31.96029466666667
$ws.Range("N5").Value = 95.88088399999999
$ws.Range("O5").Value = 0.3907265741426954
$ws.Range("P5").Value = 0.3907265741426953
$ws.Range("Q5").Value = 526.6024137014618
$ws.Range("R5").Value = 4739.421723313157
$ws.Range("S5").Value = 0.1202162571243057
$ws.Range("T5").Value = 0.1202162571243057

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.47676966666667
$ws.Range("H6").Value = 49.43030900000001
$ws.Range("I6").Value = 0.3076736139282969
$ws.Range("J6").Value = 0.3076736139282968
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 26.34807
$ws.Range("N6").Value = 79.04420999999999
$ws.Range("O6").Value = 0.3221150253382706
$ws.Range("P6").Value = 0.3221150253382706
$ws.Range("Q6").Value = 434.13108055121
$ws.Range("R6").Value = 3907.17972496089
$ws.Range("S6").Value = 0.09910629394643064
$ws.Range("T6").Value = 0.09910629394643061

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.47676966666667
$ws.Range("H7").Value = 49.43030900000001
$ws.Range("I7").Value = 0.3076736139282969
$ws.Range("J7").Value = 0.3076736139282968
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 23.48872
$ws.Range("N7").Value = 70.46616
$ws.Range("O7").Value = 0.287158400519034
$ws.Range("P7").Value = 0.287158400519034
$ws.Range("Q7").Value = 387.0182292048268
$ws.Range("R7").Value = 3483.164062843441
$ws.Range("S7").Value = 0.08835106285756052
$ws.Range("T7").Value = 0.08835106285756049

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 21.406785
$ws.Range("H8").Value = 64.220355
$ws.Range("I8").Value = 0.3997326561443945
$ws.Range("J8").Value = 0.3997326561443944
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 31.96029466666667
$ws.Range("N8").Value = 95.88088399999999
$ws.Range("O8").Value = 0.3907265741426954
$ws.Range("P8").Value = 0.3907265741426953
$ws.Range("Q8").Value = 684.1671564659799
$ws.Range("R8").Value = 6157.504408193819
$ws.Range("S8").Value = 0.1561861713082593
$ws.Range("T8").Value = 0.1561861713082593

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 21.406785
$ws.Range("H9").Value = 64.220355
$ws.Range("I9").Value = 0.3997326561443945
$ws.Range("J9").Value = 0.3997326561443944
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 26.34807
$ws.Range("N9").Value = 79.04420999999999
$ws.Range("O9").Value = 0.3221150253382706
$ws.Range("P9").Value = 0.3221150253382706
$ws.Range("Q9").Value = 564.0274696549499
$ws.Range("R9").Value = 5076.247226894549
$ws.Range("S9").Value = 0.1287598946624858
$ws.Range("T9").Value = 0.1287598946624858

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 21.406785
$ws.Range("H10").Value = 64.220355
$ws.Range("I10").Value = 0.3997326561443945
$ws.Range("J10").Value = 0.3997326561443944
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.48872
$ws.Range("N10").Value = 70.46616
$ws.Range("O10").Value = 0.287158400519034
$ws.Range("P10").Value = 0.287158400519034
$ws.Range("Q10").Value = 502.8179789652
$ws.Range("R10").Value = 4525.3618106868
$ws.Range("S10").Value = 0.1147865901736493
$ws.Range("T10").Value = 0.1147865901736493
